$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (in descending order so row indices remain valid):
# Row 13: even_MAG-GUT945.fa
# Row 11: even_MAG-GUT43457.fa
# Row 8:  even_MAG-GUT3499.fa
$ws.Rows(13).Delete()
$ws.Rows(11).Delete()
$ws.Rows(8).Delete()
